$d = $word.ActiveDocument

# Locate the "Thank you," paragraph (the last occurrence of this exact text)
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Thank you,`r") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$rng = $target.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# The newly created (still empty) paragraph now sits right after the target
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newRange = $newPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:instrText xml:space="preserve"> MERGEFIELD  Signature  \* MERGEFORMAT </w:instrText></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>«Signature»</w:t></w:r><w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p>'

$newRange.InsertXML($xml)
